$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the Periodo Mora (E) and Valor Mora (F) values between row 16 and row 18.
# Row 17 (E17="2202", F17=36341) stays unchanged.
$ws.Range("E16").Value = "2201"
$ws.Range("F16").Value = 36341

$ws.Range("E18").Value = "2203"
$ws.Range("F18").Value = 32707
